$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =====================================================================
# Phase 1: write all new cell VALUES first, in the specific order that
# reproduces the target shared-string table ordering (new unique
# strings must appear for the first time in this sequence).
# =====================================================================

# Column B - Date
$ws.Range("B9").Value = "20-12-2018"
$ws.Range("B10").Value = "21-12-2018"

# Column C - Name of employee (existing shared string, reused both rows)
$ws.Range("C9").Value = "D.Venkatesh"
$ws.Range("C10").Value = "D.Venkatesh"

# Column D - Assigned project
$ws.Range("D9").Value = "Voice Recognization System"
$ws.Range("D10").Value = "Voice Recognization System"

# Column E - Project information / learning (row 9 only for now)
$ws.Range("E9").Value = "Working On Speech to Text With Speech Recognization Tool"

# Column H / I - office log-in / log-out for row 9 (kept as text, like other rows)
$ws.Range("H9").Value = "8:30"
$ws.Range("I9").Value = "5:00"

# Row 10 remaining columns
$ws.Range("E10").Value = "Working On Speech to Text With Speech Recognization Tool & Collecting Speech of Narendra Modi "
# E10 reuses the same existing wrap/Note style as E9 (only the value differs)
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H10").Value = 0.41666666666666669
$ws.Range("I10").Value = 0.29166666666666669

$ws.Range("J10").Value = "Converting Speech to Text is done & Collected Various Speech of Narendra Modi."

# Column A - sequence number for row 10
$ws.Range("A10").Value = 7

# =====================================================================
# Phase 2: apply cell FORMATTING, in the order that reproduces the
# target style (cellXfs) table ordering. Each brand-new style should be
# produced by exactly one cell's transition; any additional cell that
# needs the identical resulting style copies the format instead of
# re-deriving it through separate property writes (which would leave a
# transient, unused style behind in the table).
# =====================================================================

# New style: wrap text (row 10 remark)
$ws.Range("J10").WrapText = $true

# New style: centered text (row 10 employee name)
$ws.Range("C10").HorizontalAlignment = -4108

# New style: time format (h:mm) + left alignment, for the log-in cell;
# the log-out cell then reuses the exact same resulting style.
$ws.Range("H10").NumberFormat = "h:mm"
$ws.Range("H10").HorizontalAlignment = -4131
$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# =====================================================================
# Phase 3: row height for the taller, wrapped row, and final selection.
# =====================================================================
$ws.Rows.Item(10).RowHeight = 38.25

$ws.Range("A10").Select()
